# Weekly update: insert a new week's record at row 12 and push the
# existing rows (12-47) down by one, so the oldest record that was at
# the bottom (row 47) becomes the new row 48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12; this shifts rows
# 12..47 down to 13..48 (Excel also extends the used range / dimension).
$ws.Rows.Item(12).Insert()

# Fill in the new row 12 with this week's data.
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value = "Maule"
$ws.Cells.Item(12, 4).Value = 44453
$ws.Cells.Item(12, 5).Value = 7
$ws.Cells.Item(12, 6).Value = 100112013
$ws.Cells.Item(12, 7).Value = "Alcachofa"
$ws.Cells.Item(12, 8).Value = "Madrigal"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 12000
$ws.Cells.Item(12, 12).Value = 12000
$ws.Cells.Item(12, 13).Value = 12000
$ws.Cells.Item(12, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 300
$ws.Cells.Item(12, 17).Value = 40
$ws.Cells.Item(12, 18).Value = "Hortaliza"
